# Rename the "Login" sheet to "info"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Name = "info"

# Add a new data row (Canada / Dina / Female) after the existing data row
$ws.Range("A3").Value = "Canada"
$ws.Range("B3").Value = "Dina"
$ws.Range("C3").Value = "Female"

# Update the active selection to match the saved view state
$ws.Range("C7").Select()
